$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elements")

# Insert a new column before the current column B ("Slice Name"), which
# becomes the new "Path" column while the existing column A becomes "ID"
# (every data row duplicates its Path value into the new ID column).
$ws.Columns("B").Insert()

# Capture the existing (pre-edit) column A values for rows 1-6 before
# overwriting the header, since row 1's new B1 should hold the old header
# text ("Path") and A1 becomes "ID".
$rowCount = 6
$idValues = @{}
for ($r = 1; $r -le $rowCount; $r++) {
    $idValues[$r] = $ws.Range("A$r").Value()
}

# Row 1 headers: A1 becomes "ID", B1 becomes "Path" (the old A1 header).
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = $idValues[1]

# Data rows 2-6: column B (new "Path" column) duplicates column A (the
# element id / path value).
for ($r = 2; $r -le $rowCount; $r++) {
    $ws.Range("B$r").Value = $idValues[$r]
}

# Fix the "Type(s)" value (now column K after the insert) for the root
# Extension element row (row 2): it should contain a single newline
# instead of being blank.
$ws.Range("K2").Value = "`n"
